# "time removed from table"
# The Start/End date columns (B, C) were date-serial values formatted as
# dates. They are converted to plain text strings ("yyyy-mm-dd", no time
# component) using a TEXT() helper formula copied in and pasted as values,
# one cell at a time so the resulting shared-string table is built up in
# the same row-by-row order Excel produced originally.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("Start") -> text dates, row by row.
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("D1").Formula = "=TEXT(B$r,""yyyy-mm-dd"")"
    $ws.Range("D1").Copy()
    $ws.Range("B$r").PasteSpecial(-4163)
}

# Column C ("End") -> text dates, row by row.
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("D1").Formula = "=TEXT(C$r,""yyyy-mm-dd"")"
    $ws.Range("D1").Copy()
    $ws.Range("C$r").PasteSpecial(-4163)
}

# Clean up the scratch helper cell used above.
$ws.Range("D1").Clear()
$excel.CutCopyMode = $false

# Leave the selection on the (now empty) helper columns D:E that were used
# while building the text values, active cell on the last-touched column.
$null = $ws.Range("D:D").Select()
$null = $ws.Range("E:E").Select()
